$d = $word.ActiveDocument

# Locate the paragraph that contains the placeholder "Abstract here!" text
# and remove the entire paragraph (including its paragraph mark) so the
# document ends right after "...State Water Project."
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "Abstract here!") {
        $para.Range.Delete()
        break
    }
}
